# Apply updated "dSF" (column F) values to Sheet1.
# Data was repulled from source; only column F values change for most rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = -1
    4  = -2
    5  = 3
    6  = -2
    7  = 1
    8  = -3
    9  = -5
    10 = -3
    11 = 2
    12 = 0
    13 = -5
    14 = 1
    15 = 3
    17 = 1
    18 = 1
    19 = -4
    20 = -1
    21 = -6
    22 = -6
    24 = 6
    25 = 4
    26 = 3
    27 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
